$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows (2..18) down to (3..19).
# This duplicates formatting (incl. the date number format) from the row
# that used to be row 2 down into the new row 3, and keeps the sheet's
# formatting intact for the rest of the table.
$ws.Rows.Item(2).Insert(-4121)  # xlShiftDown

# Full data set (after the bugfix) for rows 2..19, columns A:E
$data = @(
    @(39400, 2007, -0.7015558851707349, 2008, 1.245556224328537),
    @(39765, 2008, 0.3590181115727287, 2009, 2.957479223435766),
    @(40130, 2009, -0.01587181126743165, 2010, 3.835144370775678),
    @(40494, 2010, -0.02256889165885845, 2011, -0.6757980944263386),
    @(40862, 2011, 0.09611428386597787, 2012, -1.14533080410083),
    @(41228, 2012, -0.1827723404408288, 2013, -0.9712385417624092),
    @(41592, 2013, -0.001350220946483294, 2014, 0.6008487920565297),
    @(41957, 2014, -0.075754880139145, 2015, -1.589980933184099),
    @(42321, 2015, -0.5761528471665445, 2016, 0.4501721032283079),
    @(42689, 2016, -0.2011999787958185, 2017, -0.1501751374934579),
    @(43053, 2017, 0.1213692818849532, 2018, 1.204188598110267),
    @(43418, 2018, 0.1493219406571766, 2019, 2.372078088364726),
    @(43783, 2019, -0.4278219446121612, 2020, -2.378564786744741),
    @(44159, 2020, -1.026566979837418, 2021, -3.329887238705409),
    @(44525, 2021, 0.3179894933462268, 2022, 0.09136665742732752),
    @(44890, 2022, 0.463604920919658, 2023, -1.086119702215593),
    @(45254, 2023, 0.621639092134818, 2024, -0.185520708596465),
    @(45618, 2024, -0.6768900623516982, 2025, 0.9196327481985289)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# The row insert above copies the header row's (bold / centered) formatting
# down into the new row's B:E cells. Strip that back to the plain,
# unformatted style used by every other data row in columns B:E.
$ws.Range("B2:E2").ClearFormats() | Out-Null

# Keep the date column's formatting (number format / border / etc.)
# consistent with the rest of column A.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
